$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the units sub-header row (old row 2: "(m3/s)", "(MW)", "(MW)",
# "(GWh)", "(GWh)", "(GWh)") so the plant data rows shift up by one and
# the column headers in row 1 sit directly above the data.
$ws.Rows(2).Delete()

$ws.Range("A2:K2").Select()
